# FeatureMapCalculator.xlsx - "Added links for Course Project"
#
# The sheet computes, for a stack of Conv/MaxPool/ConvTranspose/Upsampling
# layers, the output image length from the layer's parameters
# (Padding/Stride/Kernel/Dilation/OutputPadding), chaining each layer's
# output (column B) into the next layer's input (column C).
#
# This change tweaks the Stride ("F") values of the three Conv rows in the
# second example block (rows 19/21/23) and appends a new Conv row (24) that
# continues the chain from row 23's output, reusing the already-present
# "Conv" shared string / header style used by the other "Conv" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Stride (column F) for the existing Conv rows -------------------
$ws.Range("F19").Value = 5
$ws.Range("F21").Value = 5
$ws.Range("F23").Value = 2

# --- Append a new Conv layer row at row 24 ---------------------------------
$ws.Range("A24").Value = "Conv"
$ws.Range("A24").Font.Name = $ws.Range("A19").Font.Name()
$ws.Range("A24").Font.Size = $ws.Range("A19").Font.Size()
$ws.Range("A24").Font.Bold = $ws.Range("A19").Font.Bold()

$ws.Range("B24").Formula = "=((C24+2*D24-G24*(F24-1)-1)/E24)+1"
$ws.Range("C24").Formula = "=B23"
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 1
# H24 already holds 0 from the template and is left untouched.

# --- Move the active selection to F19 (matches the saved view state) ------
[void]$ws.Range("F19").Select()

Write-Output "Applied Course Project Conv row updates"
